{"js": "// Fix small error in the ranovas (quantitative genetics) table: the QST\n// column (3rd data column, index 2) had a set of values that needed\n// correcting. Each row is located by its label (column 0) so the right\n// cell -- and only that cell -- is updated, even though some of the\n// old/new numbers are duplicated elsewhere in the table (e.g. \"0.183\",\n// \"0.010\", \"0.003\" also appear in the h2/CVA columns on other rows).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Map: row label (column 1 text) -> new QST (column 3, index 2) value\nconst updates = {\n  \"Latex exudation\": \"0.102\",\n  \"Herbivory before flowering (binary)\": \"0.025\",\n  \"Herbivory before flowering (quantitative)\": \"0.424\",\n  \"Weevil damage (quantitative)\": \"0.101\",\n  \"Flowering success\": \"0.078\",\n  \"Follicles\": \"0.011\",\n  \"Date of first follicle\": \"0.005\",\n  \"Inflorescences\": \"0.275\",\n  \"L. asclepiadis abundance\": \"0.021\",\n  \"LDMC\": \"0.110\",\n  \"SLA\": \"0.601\",\n  \"Height before flowering\": \"0.001\",\n  \"Ramets before flowering\": \"0.011\",\n  \"Mortality\": \"0.030\",\n};\n\nconst qstColIndex = 2;\nconst values = table.values;\n\nfor (let rowIndex = 0; rowIndex < values.length; rowIndex++) {\n  const label = values[rowIndex][0];\n  if (Object.prototype.hasOwnProperty.call(updates, label)) {\n    const cell = table.getCell(rowIndex, qstColIndex);\n    cell.value = updates[label];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fix small error in the ranovas (quantitative genetics) table: the QST\n# column (3rd data column) had a set of values that needed correcting.\n# Each row is located by its label in column 1 so the right cell (and\n# only that cell) is updated, even though some of the old/new numbers\n# are duplicated elsewhere in the table (e.g. \"0.183\", \"0.010\", \"0.003\").\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Map: row label (column 1 text) -> new QST (column 3) value\n$updates = @{\n    \"Latex exudation\"                           = \"0.102\"\n    \"Herbivory before flowering (binary)\"       = \"0.025\"\n    \"Herbivory before flowering (quantitative)\" = \"0.424\"\n    \"Weevil damage (quantitative)\"               = \"0.101\"\n    \"Flowering success\"                          = \"0.078\"\n    \"Follicles\"                                  = \"0.011\"\n    \"Date of first follicle\"                     = \"0.005\"\n    \"Inflorescences\"                             = \"0.275\"\n    \"L. asclepiadis abundance\"                   = \"0.021\"\n    \"LDMC\"                                       = \"0.110\"\n    \"SLA\"                                        = \"0.601\"\n    \"Height before flowering\"                    = \"0.001\"\n    \"Ramets before flowering\"                    = \"0.011\"\n    \"Mortality\"                                  = \"0.030\"\n}\n\n$qstCol = 3\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $labelCell = $t.Cell($r, 1)\n    $label = $labelCell.Range.Text\n    $label = $label -replace \"[\\u0007\\r]\", \"\"\n    if ($updates.ContainsKey($label)) {\n        $newVal = $updates[$label]\n        $t.Cell($r, $qstCol).Range.Text = $newVal\n    }\n}\n\nWrite-Output \"done\"\n"}
